$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "'60.427.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.84%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "'2.602.35"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.07%  "

# Row 4
$ws.Range("E4").Value = "  -0.16%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "'514.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.81%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "'153.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.72%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "'0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.29%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "'0.599"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.88%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "'2.614.99"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.26%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "'6.71"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.09%  "

# Row 11
$ws.Range("E11").Value = "  +0.11%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "'0.345"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.51%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "'0.129"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.89%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "'3.062.31"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.10%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "'60.487.65"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.68%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "'21.62"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.12%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "'0.0000140"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.96%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "'2.615.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.03%  "

# Row 19
$ws.Range("E19").Value = "  -0.51%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "'356.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.69%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "'10.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.56%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "'6.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.14%  "

# Row 23
$ws.Range("E23").Value = "  +0.33%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "'61.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.86%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "'0.425"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.19%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "'2.732.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.88%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "'0.166"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.57%  "

# Row 28
$ws.Range("E28").Value = "  +0.18%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "'0.0" + [char]0x2083 + "0840"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.56%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "'7.32"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.01%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "'1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.12%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "'19.42"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.93%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "'1.58"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.22%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "'150.64"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.22%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "'5.90"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.53%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "'4.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.34%  "

# Row 37
$ws.Range("E37").Value = "  -0.87%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "'0.883"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.17%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "'1.48"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.36%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "'0.843"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.17%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "'36.18"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.10%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "'3.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.87%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "'290.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.10%  "

# Row 44
$ws.Range("E44").Value = "  +1.83%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "'0.621"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.26%  "

# Row 46
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "'0.997"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.45%  "

# Row 47
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "'0.0555"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.06%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "'19.62"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.49%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "'4.95"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.63%  "

# Row 50
$ws.Range("E50").Value = "  +1.46%  "

# Row 51
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "'10.30"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.29%  "
